# PLM Iteration1 evaluation - update Man Hours and Defect-each-Task figures,
# matching the author's "update PLM Iteration1 evaluation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iteration1")

# Column E = "Man Hours": bump each feature's logged hours.
$ws.Range("E4").Value = 26
$ws.Range("E5").Value = 24
$ws.Range("E6").Value = 17
$ws.Range("E7").Value = 22
$ws.Range("E8").Value = 28
$ws.Range("E9").Value = 14

# Column I = "Defect each Task": these cells were previously blank; the
# review recorded zero defects for every task.
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 0

# Leave the cursor where the author last left it before saving.
$ws.Range("I10").Select()
